$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Outputs" sheet: the winch solenoid channel/port used to show "???"
# (unknown) for rows 6 and 7; now the channel numbers are known.
# ---------------------------------------------------------------
$wsOutputs = $wb.Worksheets.Item("Outputs")

$wsOutputs.Range("C6").Value2 = 4
$wsOutputs.Range("D6").Value2 = 4

$wsOutputs.Range("C7").Value2 = 4
$wsOutputs.Range("D7").Value2 = 3

# Update the remembered selection on that sheet.
$wsOutputs.Activate()
$wsOutputs.Range("A19").Select()

# ---------------------------------------------------------------
# "Other Inputs" sheet: the winch switch ("Do we want to turn on the
# winch motor?") now also documents which module it lives on -
# merge C8:D8 and label it "On Cyprus", centered.
# ---------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("Other Inputs")

$wsOther.Range("C8:D8").HorizontalAlignment = -4108   # xlCenter
$wsOther.Range("C8:D8").MergeCells = $true
$wsOther.Range("C8").Value2 = "On Cyprus"

# Update the remembered selection on that sheet.
$wsOther.Activate()
$wsOther.Range("C9").Select()
